# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Three country pairs swap rank (their label & stats trade rows) because
#   case totals overtook one another:
#     Malaui <-> Cabo Verde   (rows 116/117)
#     Bahamas <-> Aruba       (rows 135/136)
#     Timor Oriental <-> Santa Lucia (rows 207/208)
# - Refresh the numeric COVID figures (Casos totales, Nuevos casos,
#   Casos activos, Recuperados, Muertes hoy, Muertes) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 22:39"

# --- Country label swaps (ranking changed) -----------------------------
$ws.Range("A116").Value = "Cabo Verde"
$ws.Range("A117").Value = "Malaui"

$ws.Range("A135").Value = "Aruba"
$ws.Range("A136").Value = "Bahamas"

$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"

# --- Numeric data refresh ------------------------------------------------
$updates = @{
    "B4"   = 7316074
    "C4"   = 28513
    "D4"   = 4540053
    "E4"   = 2566609
    "G4"   = 235
    "H4"   = 209412

    "B5"   = 6073348
    "C5"   = 82767
    "D5"   = 5013367
    "E5"   = 964407
    "G5"   = 1040
    "H5"   = 95574

    "B13"  = 670766
    "C13"  = 1268
    "D13"  = 603721
    "E13"  = 50647
    "G13"  = 22
    "H13"  = 16398

    "B25"  = 286338
    "C25"  = 1313
    "E25"  = 26004

    "B53"  = 73332
    "C53"  = 632
    "D53"  = 30363
    "E53"  = 41799
    "G53"  = 5
    "H53"  = 1170

    "B116" = 5771
    "C116" = 70
    "D116" = 5031
    "E116" = 683
    "G116" = 1
    "H116" = 57

    "B117" = 5768
    "C117" = 2
    "D117" = 4206
    "E117" = 1383
    "H117" = 179

    "B119" = 5431
    "C119" = 12
    "D119" = 4821
    "E119" = 502

    "B128" = 4820
    "C128" = 9
    "D128" = 3099
    "E128" = 1692

    "B130" = 4718
    "C130" = 46
    "D130" = 1707
    "E130" = 2837
    "G130" = 3
    "H130" = 174

    "B134" = 4072
    "C134" = 34
    "D134" = 1062
    "E134" = 2818
    "G134" = 4
    "H134" = 192

    "B135" = 3844
    "C135" = 12
    "D135" = 2948
    "E135" = 871
    "H135" = 25

    "B136" = 3838
    "C136" = 48
    "D136" = 2005
    "E136" = 1744
    "H136" = 89

    "D138" = 2946
    "E138" = 543

    "B139" = 3569
    "C139" = 14
    "D139" = 2161
    "E139" = 1297
    "G139" = 1
    "H139" = 111

    "B144" = 3086
    "C144" = 6
    "D144" = 2420
    "E144" = 536

    "B167" = 1178
    "C167" = 1
    "G167" = 1
    "H167" = 84
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
